# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# described by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to stay a text value (many of these look like
    # numbers, e.g. '0.998' or '313.42') instead of letting Excel's
    # auto-detection coerce the assignment into a numeric cell, then
    # restore the 'Normal' style so no stray number-format style is
    # left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '42.856.13'
Set-TextValue $ws.Range("E2") '  +1.43%  '
Set-TextValue $ws.Range("D3") '2.282.22'
Set-TextValue $ws.Range("E3") '  -0.75%  '
Set-TextValue $ws.Range("D4") '0.998'
Set-TextValue $ws.Range("E4") '  -0.33%  '
Set-TextValue $ws.Range("D5") '313.42'
Set-TextValue $ws.Range("E5") '  -0.69%  '
Set-TextValue $ws.Range("D6") '105.18'
Set-TextValue $ws.Range("E6") '  +1.19%  '
Set-TextValue $ws.Range("E7") '  -0.68%  '
Set-TextValue $ws.Range("E8") '  +0.37%  '
Set-TextValue $ws.Range("D9") '0.601'
Set-TextValue $ws.Range("E9") '  -1.09%  '
Set-TextValue $ws.Range("D10") '39.48'
Set-TextValue $ws.Range("E10") '  -0.90%  '
Set-TextValue $ws.Range("E11") '  -0.73%  '
Set-TextValue $ws.Range("D12") '8.37'
Set-TextValue $ws.Range("E12") '  +0.52%  '
Set-TextValue $ws.Range("E13") '  +2.49%  '
Set-TextValue $ws.Range("D14") '0.992'
Set-TextValue $ws.Range("E14") '  +2.97%  '
Set-TextValue $ws.Range("D15") '15.17'
Set-TextValue $ws.Range("E15") '  -0.82%  '
Set-TextValue $ws.Range("D16") '2.629.16'
Set-TextValue $ws.Range("E16") '  -0.73%  '
Set-TextValue $ws.Range("D17") '2.290.68'
Set-TextValue $ws.Range("E17") '  -0.76%  '
Set-TextValue $ws.Range("D18") '42.581.67'
Set-TextValue $ws.Range("E18") '  +0.53%  '
Set-TextValue $ws.Range("D19") '7.35'
Set-TextValue $ws.Range("E19") '  -1.39%  '
Set-TextValue $ws.Range("E20") '  -0.60%  '
Set-TextValue $ws.Range("D21") '13.62'
Set-TextValue $ws.Range("E21") '  +22.58%  '
Set-TextValue $ws.Range("D22") '73.81'
Set-TextValue $ws.Range("E22") '  +0.74%  '
Set-TextValue $ws.Range("E23") '  +0.04%  '
Set-TextValue $ws.Range("D24") '264.07'
Set-TextValue $ws.Range("E24") '  -4.34%  '
Set-TextValue $ws.Range("D25") '2.20'
Set-TextValue $ws.Range("E25") '  -2.82%  '
Set-TextValue $ws.Range("D26") '1.00'
Set-TextValue $ws.Range("E26") '  +0.20%  '
Set-TextValue $ws.Range("D27") '10.82'
Set-TextValue $ws.Range("E27") '  +0.09%  '
Set-TextValue $ws.Range("D28") '7.08'
Set-TextValue $ws.Range("E28") '  +20.82%  '
Set-TextValue $ws.Range("D29") '2.34'
Set-TextValue $ws.Range("E29") '  -0.35%  '
Set-TextValue $ws.Range("D30") '22.43'
Set-TextValue $ws.Range("E30") '  -1.46%  '
Set-TextValue $ws.Range("D31") '36.96'
Set-TextValue $ws.Range("E31") '  +3.78%  '
Set-TextValue $ws.Range("D32") '166.22'
Set-TextValue $ws.Range("E32") '  +0.58%  '
Set-TextValue $ws.Range("D33") '0.0869'
Set-TextValue $ws.Range("E33") '  -0.12%  '
Set-TextValue $ws.Range("D34") '0.130'
Set-TextValue $ws.Range("E34") '  -3.14%  '
Set-TextValue $ws.Range("E35") '  +0.21%  '
Set-TextValue $ws.Range("D36") '0.113'
Set-TextValue $ws.Range("E36") '  -3.55%  '
Set-TextValue $ws.Range("D37") '4.52'
Set-TextValue $ws.Range("E37") '  -1.03%  '
Set-TextValue $ws.Range("E38") '  -4.41%  '
Set-TextValue $ws.Range("D39") '3.78'
Set-TextValue $ws.Range("E39") '  +1.73%  '
Set-TextValue $ws.Range("E40") '  -3.84%  '
Set-TextValue $ws.Range("E41") '  +5.23%  '
Set-TextValue $ws.Range("D42") '70.42'
Set-TextValue $ws.Range("E42") '  +1.33%  '
Set-TextValue $ws.Range("E43") '  +1.92%  '
Set-TextValue $ws.Range("D44") '94.52'
Set-TextValue $ws.Range("E44") '  -0.44%  '
Set-TextValue $ws.Range("D45") '0.999'
Set-TextValue $ws.Range("E45") '  -0.46%  '
Set-TextValue $ws.Range("D46") '12.16'
Set-TextValue $ws.Range("E46") '  +1.06%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D47") '1.736.69'
Set-TextValue $ws.Range("E47") '  +9.16%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D48") '113.28'
Set-TextValue $ws.Range("E48") '  +0.58%  '
Set-TextValue $ws.Range("D49") '79.03'
Set-TextValue $ws.Range("E49") '  -3.09%  '
Set-TextValue $ws.Range("D50") '8.71'
Set-TextValue $ws.Range("E50") '  -2.26%  '
Set-TextValue $ws.Range("D51") '5.19'
Set-TextValue $ws.Range("E51") '  +0.22%  '
